# Add a "Save" column (H) to the s_vals sheet, matching the header style
# already used by the other header cells (B1:G1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1: copy formatting from the existing "sum" header (G1),
# then set its own text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Data values for the new "Save" column (H2:H12)
$saveValues = @(0, 0, 0, 0, 0, 1, 0, 0, 0, 0, 1)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
